$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ E=3; G=39.27015933333333; H=117.810478; I=0.2257020976862494; J=0.2257020976862494; K=3; M=103.4766596666667; N=310.429979; O=0.877785331764719; P=0.8777853317647188; Q=4063.544912391107; R=36571.90421151996; S=0.1981179906975175; T=0.1981179906975175 }
    3  = @{ E=3; G=39.27015933333333; H=117.810478; I=0.2257020976862494; J=0.2257020976862494; K=3; M=0.8265796666666668; N=2.479739; O=0.007011818020336602; P=0.0070118180203366; Q=32.45991521169356; R=292.139236905242; S=0.001582582035784216; T=0.001582582035784215 }
    4  = @{ E=3; G=39.27015933333333; H=117.810478; I=0.2257020976862494; J=0.2257020976862494; K=3; M=13.58054833333333; N=40.741645; O=0.1152028502149446; P=0.1152028502149446; Q=533.3102968840344; R=4799.792671956309; S=0.02600152495294779; T=0.02600152495294778 }
    5  = @{ E=3; G=119.3024773333333; H=357.907432; I=0.6856814398113102; J=0.6856814398113102; K=3; M=103.4766596666667; N=310.429979; O=0.877785331764719; P=0.8777853317647188; Q=12345.02184441155; R=111105.1965997039; S=0.6018811101296812; T=0.601881110129681 }
    6  = @{ E=3; G=119.3024773333333; H=357.907432; I=0.6856814398113102; J=0.6856814398113102; K=3; M=0.8265796666666668; N=2.479739; O=0.007011818020336602; P=0.0070118180203366; Q=98.61300194669425; R=887.5170175202481; S=0.004807873475879292; T=0.004807873475879291 }
    7  = @{ E=3; G=119.3024773333333; H=357.907432; I=0.6856814398113102; J=0.6856814398113102; K=3; M=13.58054833333333; N=40.741645; O=0.1152028502149446; P=0.1152028502149446; Q=1620.193059711738; R=14581.73753740564; S=0.07899245620574993; T=0.07899245620574992 }
    8  = @{ E=3; G=15.418477; H=46.255431; I=0.08861646250244033; J=0.08861646250244033; K=3; M=103.4766596666667; N=310.429979; O=0.877785331764719; P=0.8777853317647188; Q=1595.452497107328; R=14359.07247396595; S=0.07778623093752036; T=0.07778623093752035 }
    9  = @{ E=3; G=15.418477; H=46.255431; I=0.08861646250244033; J=0.08861646250244033; K=3; M=0.8265796666666668; N=2.479739; O=0.007011818020336602; P=0.0070118180203366; Q=12.74459957916767; R=114.701396212509; S=0.0006213625086730939; T=0.0006213625086730938 }
    10 = @{ E=3; G=15.418477; H=46.255431; I=0.08861646250244033; J=0.08861646250244033; K=3; M=13.58054833333333; N=40.741645; O=0.1152028502149446; P=0.1152028502149446; Q=209.3913721248883; R=1884.522349123995; S=0.01020886905624689; T=0.01020886905624689 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
